$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header: "Sub Categoria" -> "SubCategoria"
$ws.Range("B1").Value = "SubCategoria"

# Replace placeholder "_" with new category name
$ws.Range("A2").Value = "Incidencias / Errores"

# Remove the now-redundant duplicate last row (row 36 duplicated row 35's value)
$ws.Rows.Item(36).Delete()

# Update the saved selection to match the target state
$ws.Range("E30").Select()
